$wb = $excel.ActiveWorkbook
Write-Host ("Path: " + $wb.Path())
Write-Host ("FullName: " + $wb.FullName())
Write-Host ("Name: " + $wb.Name())
